$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.828.38'
$ws.Range('E2').Value = '  +0.34%  '
$ws.Range('D3').Value = '2.539.42'
$ws.Range('E3').Value = '  +0.09%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '304.27'
$ws.Range('D6').Value = '98.52'
$ws.Range('E6').Value = '  +5.45%  '
$ws.Range('D7').Value = '0.577'
$ws.Range('E7').Value = '  +0.81%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '0.545'
$ws.Range('E9').Value = '  +0.10%  '
$ws.Range('D10').Value = '36.52'
$ws.Range('E10').Value = '  +1.85%  '
$ws.Range('E11').Value = '  +2.72%  '
$ws.Range('D12').Value = '7.64'
$ws.Range('E12').Value = '  -0.55%  '
$ws.Range('E13').Value = '  -1.03%  '
$ws.Range('D14').Value = '2.931.79'
$ws.Range('E14').Value = '  +0.13%  '
$ws.Range('D15').Value = '2.545.16'
$ws.Range('E15').Value = '  +0.01%  '
$ws.Range('D16').Value = '15.20'
$ws.Range('E16').Value = '  +7.41%  '
$ws.Range('E17').Value = '  +0.02%  '
$ws.Range('D18').Value = '42.870.12'
$ws.Range('E18').Value = '  +0.31%  '
$ws.Range('D19').Value = '13.19'
$ws.Range('E19').Value = '  +4.91%  '
$ws.Range('D20').Value = '0.0₃0989'
$ws.Range('E20').Value = '  +1.27%  '
$ws.Range('E21').Value = '  +0.75%  '
$ws.Range('D22').Value = '71.72'
$ws.Range('E22').Value = '  +0.15%  '
$ws.Range('D23').Value = '253.69'
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('E24').Value = '  +1.16%  '
$ws.Range('D25').Value = '2.06'
$ws.Range('E25').Value = '  -2.63%  '
$ws.Range('D26').Value = '27.75'
$ws.Range('E26').Value = '  -4.12%  '
$ws.Range('E27').Value = '  -0.25%  '
$ws.Range('E28').Value = '  +9.31%  '
$ws.Range('E29').Value = '  +0.18%  '
$ws.Range('D30').Value = '38.61'
$ws.Range('E30').Value = '  +5.86%  '
$ws.Range('D31').Value = '6.19'
$ws.Range('E31').Value = '  +2.28%  '
$ws.Range('D32').Value = '157.18'
$ws.Range('E32').Value = '  +3.24%  '
$ws.Range('E33').Value = '  +0.19%  '
$ws.Range('D34').Value = '19.17'
$ws.Range('E34').Value = '  +11.01%  '
$ws.Range('E35').Value = '  +1.08%  '
$ws.Range('D36').Value = '3.29'
$ws.Range('E36').Value = '  -1.79%  '
$ws.Range('D37').Value = '2.63'
$ws.Range('E37').Value = '  -4.25%  '
$ws.Range('D38').Value = '0.115'
$ws.Range('E38').Value = '  +2.10%  '
$ws.Range('D39').Value = '24.94'
$ws.Range('E39').Value = '  +8.58%  '
$ws.Range('E40').Value = '  +0.81%  '
$ws.Range('E41').Value = '  +10.61%  '
$ws.Range('E42').Value = '  +1.89%  '
$ws.Range('E43').Value = '  +1.93%  '
$ws.Range('D44').Value = '2.082.08'
$ws.Range('E44').Value = '  -0.12%  '
$ws.Range('D45').Value = '0.0304'
$ws.Range('E45').Value = '  -1.09%  '
$ws.Range('D46').Value = '0.998'
$ws.Range('E46').Value = '  -0.09%  '
$ws.Range('D47').Value = '86.39'
$ws.Range('E47').Value = '  +2.88%  '
$ws.Range('D48').Value = '8.97'
$ws.Range('E48').Value = '  -0.98%  '
$ws.Range('D49').Value = '2.789.12'
$ws.Range('E49').Value = '  +0.14%  '
$ws.Range('E50').Value = '  +7.36%  '
$ws.Range('E51').Value = '  +2.34%  '
